# Update Amelx-Cd63 LR-pairs sheet with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.069782333333333
$ws.Cells.Item(2, 8).Value = 3.209347
$ws.Cells.Item(2, 9).Value = 0.5403049320348338
$ws.Cells.Item(2, 10).Value = 0.5403049320348337
$ws.Cells.Item(2, 13).Value = 9.873811666666667
$ws.Cells.Item(2, 14).Value = 29.621435
$ws.Cells.Item(2, 15).Value = 0.01897536961063408
$ws.Cells.Item(2, 16).Value = 0.01897536961063408
$ws.Cells.Item(2, 17).Value = 10.56282928366056
$ws.Cells.Item(2, 18).Value = 95.065463552945
$ws.Cells.Item(2, 19).Value = 0.0102524857878095
$ws.Cells.Item(2, 20).Value = 0.01025248578780949

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.069782333333333
$ws.Cells.Item(3, 8).Value = 3.209347
$ws.Cells.Item(3, 9).Value = 0.5403049320348338
$ws.Cells.Item(3, 10).Value = 0.5403049320348337
$ws.Cells.Item(3, 15).Value = 0.368560155467396
$ws.Cells.Item(3, 16).Value = 0.368560155467396
$ws.Cells.Item(3, 17).Value = 205.1626968456931
$ws.Cells.Item(3, 18).Value = 1846.464271611238
$ws.Cells.Item(3, 19).Value = 0.1991348697505592
$ws.Cells.Item(3, 20).Value = 0.1991348697505591

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.069782333333333
$ws.Cells.Item(4, 8).Value = 3.209347
$ws.Cells.Item(4, 9).Value = 0.5403049320348338
$ws.Cells.Item(4, 10).Value = 0.5403049320348337
$ws.Cells.Item(4, 13).Value = 140.35703
$ws.Cells.Item(4, 14).Value = 421.07109
$ws.Cells.Item(4, 15).Value = 0.2697364109842271
$ws.Cells.Item(4, 16).Value = 0.2697364109842271
$ws.Cells.Item(4, 17).Value = 150.1514710531367
$ws.Cells.Item(4, 18).Value = 1351.36323947823
$ws.Cells.Item(4, 19).Value = 0.1457399132041528
$ws.Cells.Item(4, 20).Value = 0.1457399132041528

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.069782333333333
$ws.Cells.Item(5, 8).Value = 3.209347
$ws.Cells.Item(5, 9).Value = 0.5403049320348338
$ws.Cells.Item(5, 10).Value = 0.5403049320348337
$ws.Cells.Item(5, 13).Value = 178.3381523333333
$ws.Cells.Item(5, 14).Value = 535.014457
$ws.Cells.Item(5, 15).Value = 0.3427280639377429
$ws.Cells.Item(5, 16).Value = 0.3427280639377429
$ws.Cells.Item(5, 17).Value = 190.7830047255088
$ws.Cells.Item(5, 18).Value = 1717.047042529579
$ws.Cells.Item(5, 19).Value = 0.1851776632923123
$ws.Cells.Item(5, 20).Value = 0.1851776632923123

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 7).Value = 0.663689
$ws.Cells.Item(6, 8).Value = 1.991067
$ws.Cells.Item(6, 9).Value = 0.3352031799963669
$ws.Cells.Item(6, 10).Value = 0.3352031799963669
$ws.Cells.Item(6, 13).Value = 9.873811666666667
$ws.Cells.Item(6, 14).Value = 29.621435
$ws.Cells.Item(6, 15).Value = 0.01897536961063408
$ws.Cells.Item(6, 16).Value = 0.01897536961063408
$ws.Cells.Item(6, 17).Value = 6.553140191238333
$ws.Cells.Item(6, 18).Value = 58.978261721145
$ws.Cells.Item(6, 19).Value = 0.006360604235090966
$ws.Cells.Item(6, 20).Value = 0.006360604235090965

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 7).Value = 0.663689
$ws.Cells.Item(7, 8).Value = 1.991067
$ws.Cells.Item(7, 9).Value = 0.3352031799963669
$ws.Cells.Item(7, 10).Value = 0.3352031799963669
$ws.Cells.Item(7, 15).Value = 0.368560155467396
$ws.Cells.Item(7, 16).Value = 0.368560155467396
$ws.Cells.Item(7, 17).Value = 127.2821777515687
$ws.Cells.Item(7, 18).Value = 1145.539599764118
$ws.Cells.Item(7, 19).Value = 0.1235425361326265
$ws.Cells.Item(7, 20).Value = 0.1235425361326265

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 7).Value = 0.663689
$ws.Cells.Item(8, 8).Value = 1.991067
$ws.Cells.Item(8, 9).Value = 0.3352031799963669
$ws.Cells.Item(8, 10).Value = 0.3352031799963669
$ws.Cells.Item(8, 13).Value = 140.35703
$ws.Cells.Item(8, 14).Value = 421.07109
$ws.Cells.Item(8, 15).Value = 0.2697364109842271
$ws.Cells.Item(8, 16).Value = 0.2697364109842271
$ws.Cells.Item(8, 17).Value = 93.15341688367
$ws.Cells.Item(8, 18).Value = 838.38075195303
$ws.Cells.Item(8, 19).Value = 0.09041650272271987
$ws.Cells.Item(8, 20).Value = 0.09041650272271987

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 7).Value = 0.663689
$ws.Cells.Item(9, 8).Value = 1.991067
$ws.Cells.Item(9, 9).Value = 0.3352031799963669
$ws.Cells.Item(9, 10).Value = 0.3352031799963669
$ws.Cells.Item(9, 13).Value = 178.3381523333333
$ws.Cells.Item(9, 14).Value = 535.014457
$ws.Cells.Item(9, 15).Value = 0.3427280639377429
$ws.Cells.Item(9, 16).Value = 0.3427280639377429
$ws.Cells.Item(9, 17).Value = 118.3610699839577
$ws.Cells.Item(9, 18).Value = 1065.249629855619
$ws.Cells.Item(9, 19).Value = 0.1148835369059296
$ws.Cells.Item(9, 20).Value = 0.1148835369059296

# Row 10
$ws.Cells.Item(10, 1).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.246489
$ws.Cells.Item(10, 8).Value = 0.7394670000000001
$ws.Cells.Item(10, 9).Value = 0.1244918879687994
$ws.Cells.Item(10, 10).Value = 0.1244918879687994
$ws.Cells.Item(10, 13).Value = 9.873811666666667
$ws.Cells.Item(10, 14).Value = 29.621435
$ws.Cells.Item(10, 15).Value = 0.01897536961063408
$ws.Cells.Item(10, 16).Value = 0.01897536961063408
$ws.Cells.Item(10, 17).Value = 2.433785963905001
$ws.Cells.Item(10, 18).Value = 21.904073675145
$ws.Cells.Item(10, 19).Value = 0.002362279587733618
$ws.Cells.Item(10, 20).Value = 0.002362279587733618

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.246489
$ws.Cells.Item(11, 8).Value = 0.7394670000000001
$ws.Cells.Item(11, 9).Value = 0.1244918879687994
$ws.Cells.Item(11, 10).Value = 0.1244918879687994
$ws.Cells.Item(11, 15).Value = 0.368560155467396
$ws.Cells.Item(11, 16).Value = 0.368560155467396
$ws.Cells.Item(11, 17).Value = 47.27162377530201
$ws.Cells.Item(11, 18).Value = 425.4446139777181
$ws.Cells.Item(11, 19).Value = 0.04588274958421035
$ws.Cells.Item(11, 20).Value = 0.04588274958421035

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.246489
$ws.Cells.Item(12, 8).Value = 0.7394670000000001
$ws.Cells.Item(12, 9).Value = 0.1244918879687994
$ws.Cells.Item(12, 10).Value = 0.1244918879687994
$ws.Cells.Item(12, 13).Value = 140.35703
$ws.Cells.Item(12, 14).Value = 421.07109
$ws.Cells.Item(12, 15).Value = 0.2697364109842271
$ws.Cells.Item(12, 16).Value = 0.2697364109842271
$ws.Cells.Item(12, 17).Value = 34.59646396767
$ws.Cells.Item(12, 18).Value = 311.3681757090301
$ws.Cells.Item(12, 19).Value = 0.03357999505735443
$ws.Cells.Item(12, 20).Value = 0.03357999505735442

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.246489
$ws.Cells.Item(13, 8).Value = 0.7394670000000001
$ws.Cells.Item(13, 9).Value = 0.1244918879687994
$ws.Cells.Item(13, 10).Value = 0.1244918879687994
$ws.Cells.Item(13, 13).Value = 178.3381523333333
$ws.Cells.Item(13, 14).Value = 535.014457
$ws.Cells.Item(13, 15).Value = 0.3427280639377429
$ws.Cells.Item(13, 16).Value = 0.3427280639377429
$ws.Cells.Item(13, 17).Value = 43.95839283049101
$ws.Cells.Item(13, 18).Value = 395.6255354744191
$ws.Cells.Item(13, 19).Value = 0.042666863739501
$ws.Cells.Item(13, 20).Value = 0.04266686373950099
